$d = $word.ActiveDocument

# --- 1. Replace the "Building Height" diagram picture (first InlineShape)
#        with a hyperlink to the source image.
$shape1 = $d.InlineShapes(1)
$ip1 = $d.Range($shape1.Range.Start, $shape1.Range.Start)
$shape1.Delete()
$url1 = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/CCI01_Building_Height.jpg?h=100%25&w=100%25"
$d.Hyperlinks.Add($ip1, $url1, "", "", $url1) | Out-Null

# --- 2. Replace the "Additional Height for Predominant Sky Terrace Storey"
#        diagram picture (remaining InlineShape) with a hyperlink to the
#        source image.
$shape2 = $d.InlineShapes(1)
$ip2 = $d.Range($shape2.Range.Start, $shape2.Range.Start)
$shape2.Delete()
$url2 = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C04_Additional_Height_for_Sky_Terrace_Floors.jpg?h=100%25&w=100%25"
$d.Hyperlinks.Add($ip2, $url2, "", "", $url2) | Out-Null

Write-Output "done"
